$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1488.6666
$ws.Range("J43").Value = 1485.4286
$ws.Range("L43").Value = 1485.4286
$ws.Range("N43").Value = -1623.4286
$ws.Range("H98").Value = 4520.7085
$ws.Range("I98").Value = 4172.6113
$ws.Range("K98").Value = 4172.6113
$ws.Range("M98").Value = -2674.6113
$ws.Range("H122").Value = 4520.7085
$ws.Range("I122").Value = 4172.6113
$ws.Range("K122").Value = 12517.8339
$ws.Range("M122").Value = -10067.8339
$ws.Range("H124").Value = 68988.2
$ws.Range("J124").Value = 68988.2
$ws.Range("L124").Value = 68988.2
$ws.Range("N124").Value = -78808.2
$ws.Range("H125").Value = 564.125
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 40000
$ws.Range("J126").Value = 40000
$ws.Range("L126").Value = 40000
$ws.Range("N126").Value = -49880
$ws.Range("H132").Value = 1158.7826
$ws.Range("I132").Value = 1062.05
$ws.Range("K132").Value = 3186.15
$ws.Range("M132").Value = -656.1499999999996
$ws.Range("H137").Value = 41166.2
$ws.Range("I137").Value = 789.7692
$ws.Range("K137").Value = 2369.3076
$ws.Range("M137").Value = 180.6923999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 773.8570999999999
$ws.Range("I4").Value = 741.75
$ws.Range("J4").Value = 816.6667
$ws.Range("K4").Value = 741.75
$ws.Range("L4").Value = 816.6667
$ws.Range("M4").Value = -625.75
$ws.Range("N4").Value = -1048.6667
$ws.Range("H41").Value = 6704.4
$ws.Range("I41").Value = 6704.4
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 6704.4
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -6290.4
$ws.Range("N41").ClearContents()
$ws.Range("H61").Value = 38158.043
$ws.Range("I61").Value = 48272.707
$ws.Range("K61").Value = 48272.707
$ws.Range("M61").Value = -48060.707
$ws.Range("H74").Value = 764.91895
$ws.Range("I74").Value = 552.84375
$ws.Range("K74").Value = 552.84375
$ws.Range("M74").Value = 321.15625
$ws.Range("H77").Value = 764.91895
$ws.Range("I77").Value = 552.84375
$ws.Range("K77").Value = 2764.21875
$ws.Range("M77").Value = 1603.78125
$ws.Range("H122").Value = 1207.5
$ws.Range("I122").Value = 1224.3077
$ws.Range("K122").Value = 3672.9231
$ws.Range("M122").Value = -1222.9231
$ws.Range("H136").Value = 38158.043
$ws.Range("I136").Value = 48272.707
$ws.Range("K136").Value = 144818.121
$ws.Range("M136").Value = -142268.121

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 2000
$ws.Range("I8").Value = 2000
$ws.Range("K8").Value = 2000
$ws.Range("M8").Value = -1860
$ws.Range("H19").Value = 5099.1
$ws.Range("I19").Value = 4999
$ws.Range("J19").Value = 6000
$ws.Range("K19").Value = 4999
$ws.Range("L19").Value = 6000
$ws.Range("M19").Value = -4826
$ws.Range("N19").Value = -6346
$ws.Range("H80").Value = 6492
$ws.Range("I80").Value = 85
$ws.Range("J80").Value = 8322.571
$ws.Range("K80").Value = 85
$ws.Range("L80").Value = 8322.571
$ws.Range("M80").Value = 913
$ws.Range("N80").Value = -10318.571
$ws.Range("H83").Value = 6492
$ws.Range("I83").Value = 85
$ws.Range("J83").Value = 8322.571
$ws.Range("K83").Value = 425
$ws.Range("L83").Value = 41612.855
$ws.Range("M83").Value = 4567
$ws.Range("N83").Value = -51596.855
$ws.Range("H134").Value = 7550.125
$ws.Range("I134").Value = 9377.4375
$ws.Range("J134").Value = 3895.5
$ws.Range("K134").Value = 28132.3125
$ws.Range("L134").Value = 11686.5
$ws.Range("M134").Value = -25597.3125
$ws.Range("N134").Value = -16756.5
$ws.Range("H135").Value = 59325
$ws.Range("J135").Value = 59325
$ws.Range("L135").Value = 59325
$ws.Range("N135").Value = -69465

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1399.9
$ws.Range("I22").Value = 612.75
$ws.Range("J22").Value = 1924.6666
$ws.Range("K22").Value = 612.75
$ws.Range("L22").Value = 1924.6666
$ws.Range("M22").Value = -262.75
$ws.Range("N22").Value = -2624.6666
$ws.Range("H31").Value = 2230.2964
$ws.Range("I31").Value = 1856.0625
$ws.Range("J31").Value = 2774.6365
$ws.Range("K31").Value = 1856.0625
$ws.Range("L31").Value = 2774.6365
$ws.Range("M31").Value = -1561.0625
$ws.Range("N31").Value = -3364.6365
$ws.Range("H34").Value = 2230.2964
$ws.Range("I34").Value = 1856.0625
$ws.Range("J34").Value = 2774.6365
$ws.Range("K34").Value = 1856.0625
$ws.Range("L34").Value = 2774.6365
$ws.Range("M34").Value = -1654.0625
$ws.Range("N34").Value = -3178.6365
$ws.Range("H99").Value = 2692.35
$ws.Range("J99").Value = 2717.3333
$ws.Range("L99").Value = 2717.3333
$ws.Range("N99").Value = -5713.3333
$ws.Range("H103").Value = 17497.5
$ws.Range("I103").Value = 17497.5
$ws.Range("K103").Value = 17497.5
$ws.Range("M103").Value = -16325.5
$ws.Range("H126").Value = 2692.35
$ws.Range("J126").Value = 2717.3333
$ws.Range("L126").Value = 8151.999899999999
$ws.Range("N126").Value = -13091.9999
$ws.Range("H132").Value = 1928.4375
$ws.Range("I132").Value = 1322.9
$ws.Range("K132").Value = 3968.7
$ws.Range("M132").Value = -1438.7
$ws.Range("H134").Value = 3989.889
$ws.Range("I134").Value = 3535
$ws.Range("J134").Value = 4899.6665
$ws.Range("K134").Value = 10605
$ws.Range("L134").Value = 14698.9995
$ws.Range("M134").Value = -8070
$ws.Range("N134").Value = -19768.9995

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 870
$ws.Range("I61").Value = 1297.5
$ws.Range("J61").Value = 625.7143
$ws.Range("K61").Value = 3892.5
$ws.Range("L61").Value = 1877.1429
$ws.Range("M61").Value = -3677.5
$ws.Range("N61").Value = -2307.1429
$ws.Range("H131").Value = 16843.553
$ws.Range("J131").Value = 17774.887
$ws.Range("L131").Value = 53324.66099999999
$ws.Range("N131").Value = -63404.66099999999
$ws.Range("H133").Value = 3474.875
$ws.Range("I133").Value = 2499.75
$ws.Range("J133").Value = 4450
$ws.Range("K133").Value = 7499.25
$ws.Range("L133").Value = 13350
$ws.Range("M133").Value = -2439.25
$ws.Range("N133").Value = -23470
$ws.Range("H141").Value = 2788.739
$ws.Range("I141").Value = 2605.4
$ws.Range("K141").Value = 7816.200000000001
$ws.Range("M141").Value = -2636.200000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 122.333336
$ws.Range("I2").Value = 30
$ws.Range("K2").Value = 30
$ws.Range("M2").Value = 83
$ws.Range("H46").Value = 17273.234
$ws.Range("J46").Value = 17102.812
$ws.Range("L46").Value = 17102.812
$ws.Range("N46").Value = -17414.812

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3054.5
$ws.Range("I7").Value = 3481.5
$ws.Range("K7").Value = 3481.5
$ws.Range("M7").Value = -3369.5
$ws.Range("H46").Value = 2414.2856
$ws.Range("I46").Value = 1755.375
$ws.Range("J46").Value = 3292.8333
$ws.Range("K46").Value = 1755.375
$ws.Range("L46").Value = 3292.8333
$ws.Range("M46").Value = -1567.375
$ws.Range("N46").Value = -3668.8333
$ws.Range("H93").Value = 1175.3334
$ws.Range("I93").Value = 639.8333
$ws.Range("K93").Value = 639.8333
$ws.Range("M93").Value = 608.1667
$ws.Range("H126").Value = 3054.5
$ws.Range("I126").Value = 3481.5
$ws.Range("K126").Value = 10444.5
$ws.Range("M126").Value = -7974.5
$ws.Range("H132").Value = 2295.5122
$ws.Range("I132").Value = 1807.4706
$ws.Range("K132").Value = 5422.4118
$ws.Range("M132").Value = -2892.4118
$ws.Range("H133").Value = 69500
$ws.Range("J133").Value = 69500
$ws.Range("L133").Value = 69500
$ws.Range("N133").Value = -74560

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1535.75
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

Write-Host "Applied all Tonberry Profits updates"